# "fixing lca linking to database"
#
# The "district heating - natural gas-fired boiler" (T25 / DH) row had been
# entered on the HEATING sheet, but it is actually a district-cooling
# technology (natural-gas boiler feeding an absorption chiller) and belongs
# on the COOLING sheet instead. Move that row over, correct its
# description/code/source text, and update the active sheet/selection to
# match.

$wb = $excel.ActiveWorkbook
$wsHeating = $wb.Worksheets.Item("HEATING")
$wsCooling = $wb.Worksheets.Item("COOLING")

# 1) Copy the misplaced row (with all formatting) from HEATING row 5 to the
#    first free row on COOLING (row 7).
$wsHeating.Range("A5:I5").Copy($wsCooling.Range("A7:I7"))

# 2) Remove the row from HEATING - it doesn't belong there.
$wsHeating.Rows.Item(5).Delete()

# 3) Correct the description/code/source on the newly-added COOLING row so
#    it correctly links to the district-cooling database entry.
$wsCooling.Cells.Item(7, 1).Value = "district cooling - natural gas-fired boiler for absorption chiller"
$wsCooling.Cells.Item(7, 2).Value = "T25"
$wsCooling.Cells.Item(7, 3).Value = "NG"

# 4) Update sheet selections/active tab to reflect the edits made above.
$wsHeating.Activate()
$wsHeating.Range("A5:XFD5").Select()

$wsCooling.Activate()
$wsCooling.Range("D26").Select()
